# v1.0.6: bug fixes and new feature: save last submission
#
# The spreadsheet is repurposed from tracking "linh kien" (components) to
# tracking "thanh pham" (finished goods): the three sheets are renamed and
# their header rows are updated to the new column set. All sample/demo
# data rows are cleared out (templates reset), and the now-unused trailing
# columns are dropped.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # nhap-linhkien  -> nhap-thanhpham
$ws2 = $wb.Worksheets.Item(2)   # xuat-linhkien  -> xuat-thanhpham
$ws3 = $wb.Worksheets.Item(3)   # ton-linhkien   -> ton-thanhpham

# --- Sheet 1: nhap-thanhpham ------------------------------------------------
# New 6-column header, sample rows 2:3 removed, old columns G:J dropped.
$ws1.Rows("2:3").Delete()
$ws1.Columns("G:J").Delete()
$ws1.Range("A1").Value = "Tên Hàng"
$ws1.Range("B1").Value = "MCU"
$ws1.Range("C1").Value = "Sổ Hợp Đồng"
$ws1.Range("D1").Value = "Chip"
$ws1.Range("E1").Value = "Ngày Nhập"
$ws1.Range("F1").Value = "Số Lượng"

# --- Sheet 2: xuat-thanhpham -------------------------------------------------
# Same new 6-column header; this sheet never had sample rows, just columns G:J.
$ws2.Columns("G:J").Delete()
$ws2.Range("A1").Value = "Tên Hàng"
$ws2.Range("B1").Value = "MCU"
$ws2.Range("C1").Value = "Sổ Hợp Đồng"
$ws2.Range("D1").Value = "Chip"
$ws2.Range("E1").Value = "Ngày Nhập"
$ws2.Range("F1").Value = "Số Lượng"

# --- Sheet 3: ton-thanhpham ---------------------------------------------------
# Header (Tên Hàng / Số Lượng / Đơn Vị Tính) is unchanged; just drop the
# two sample rows.
$ws3.Rows("2:3").Delete()

# --- Rename the sheets to match the new "thanh pham" naming -----------------
$ws1.Name = "nhap-thanhpham"
$ws2.Name = "xuat-thanhpham"
$ws3.Name = "ton-thanhpham"
